$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Change 1: Latest HO Xliff Generate Date on Overview sheet (G2, G3)
$wsOverview.Range("G2").Value = "2016-08-05 02:33:03"
$wsOverview.Range("G3").Value = "2016-08-05 02:33:03"

# Change 2: Latest Handoff Datetime on zh-cn sheet (H2, H3)
$wsZhCn.Range("H2").Value = "2016-08-05 02:32:50"
$wsZhCn.Range("H3").Value = "2016-08-05 02:32:50"

# Change 3 & 4: Error Detail column (P2, P3) on both zh-cn and de-de sheets.
# The commit hash embedded in the "latest" link changes from
# f8a0a8be1ae2b20593b890a0c3af2d8de672c270 to ea390eaa40dc00dc8afc058ae26070a94ecf6d1c
$errorDetailA = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/b342cccd51e418fd787d5d107c1c5f8de858cc03/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/ea390eaa40dc00dc8afc058ae26070a94ecf6d1c/e2e/a.md."
$errorDetailB = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/b342cccd51e418fd787d5d107c1c5f8de858cc03/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/ea390eaa40dc00dc8afc058ae26070a94ecf6d1c/e2e/b.md."

$wsZhCn.Range("P2").Value = $errorDetailA
$wsZhCn.Range("P3").Value = $errorDetailB

$wsDeDe.Range("P2").Value = $errorDetailA
$wsDeDe.Range("P3").Value = $errorDetailB
